# Update the "cryptos" sheet with refreshed price/volume figures (and a
# few re-sorted rows) as captured by the scheduled GitHub Actions scrape.
#
# The Price/Volume columns hold plain text (e.g. "68.963.17", "0.999",
# "  +0.32%  ") rather than numbers, so we force the cell's number format
# to Text ("@") before assigning the value -- otherwise Excel would helpfully
# "interpret" strings like "0.999" or "8.18" as real numbers. We reset the
# style back to "Normal" right after so we don't leave stray text-format
# styling on cells that didn't have any before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "68.963.17"
Set-TextValue 2 5 "  +0.32%  "

Set-TextValue 3 4 "3.745.00"
Set-TextValue 3 5 "  +0.10%  "

Set-TextValue 4 4 "0.999"
Set-TextValue 4 5 "  -0.04%  "

Set-TextValue 5 4 "600.56"
Set-TextValue 5 5 "  -0.19%  "

Set-TextValue 6 4 "165.66"
Set-TextValue 6 5 "  -1.74%  "

Set-TextValue 7 4 "3.743.33"
Set-TextValue 7 5 "  +0.06%  "

Set-TextValue 8 5 "  -0.05%  "

Set-TextValue 9 5 "  +0.92%  "

Set-TextValue 10 4 "0.172"
Set-TextValue 10 5 "  +4.97%  "

Set-TextValue 11 4 "6.41"
Set-TextValue 11 5 "  +1.11%  "

Set-TextValue 12 4 "0.459"
Set-TextValue 12 5 "  -0.53%  "

Set-TextValue 13 4 "37.73"
Set-TextValue 13 5 "  -1.18%  "

Set-TextValue 14 4 "0.0000248"
Set-TextValue 14 5 "  +0.91%  "

Set-TextValue 15 4 "4.370.02"
Set-TextValue 15 5 "  +0.10%  "

Set-TextValue 16 4 "3.747.01"
Set-TextValue 16 5 "  +0.30%  "

Set-TextValue 17 4 "69.058.13"
Set-TextValue 17 5 "  +0.46%  "

Set-TextValue 18 4 "7.43"
Set-TextValue 18 5 "  +1.76%  "

Set-TextValue 19 5 "  +3.18%  "

Set-TextValue 20 5 "  -0.84%  "

Set-TextValue 21 4 "11.27"
Set-TextValue 21 5 "  +4.99%  "

Set-TextValue 22 4 "492.16"
Set-TextValue 22 5 "  -0.47%  "

Set-TextValue 23 4 "0.725"
Set-TextValue 23 5 "  -0.38%  "

Set-TextValue 24 2 "PEPE"
Set-TextValue 24 3 "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue 24 4 "0.0000149"
Set-TextValue 24 5 "  +3.24%  "

Set-TextValue 25 2 "Litecoin"
Set-TextValue 25 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue 25 4 "84.60"
Set-TextValue 25 5 "  -0.79%  "

Set-TextValue 26 5 "  -2.02%  "

Set-TextValue 27 4 "12.33"
Set-TextValue 27 5 "  -0.80%  "

Set-TextValue 28 5 "  -0.99%  "

Set-TextValue 30 5 "  -0.49%  "

Set-TextValue 31 4 "8.18"
Set-TextValue 31 5 "  +3.32%  "

Set-TextValue 32 5 "  -4.81%  "

Set-TextValue 33 4 "31.67"
Set-TextValue 33 5 "  -0.13%  "

Set-TextValue 34 4 "3.887.20"
Set-TextValue 34 5 "  +0.05%  "

Set-TextValue 35 4 "3.679.56"
Set-TextValue 35 5 "  +0.14%  "

Set-TextValue 36 5 "  -0.30%  "

Set-TextValue 37 4 "5.94"
Set-TextValue 37 5 "  +1.98%  "

Set-TextValue 38 5 "  -0.03%  "

Set-TextValue 39 5 "  +3.94%  "

Set-TextValue 40 5 "  +0.03%  "

Set-TextValue 41 4 "3.13"
Set-TextValue 41 5 "  +9.18%  "

Set-TextValue 42 4 "0.325"
Set-TextValue 42 5 "  -0.46%  "

Set-TextValue 43 2 "Bittensor"
Set-TextValue 43 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue 43 4 "427.96"
Set-TextValue 43 5 "  -2.28%  "

Set-TextValue 44 2 "OKB"
Set-TextValue 44 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 44 4 "48.63"
Set-TextValue 44 5 "  -0.59%  "

Set-TextValue 45 5 "  +0.17%  "

Set-TextValue 46 4 "8.44"
Set-TextValue 46 5 "  -0.74%  "

Set-TextValue 47 5 "  -0.01%  "

Set-TextValue 48 2 "Arweave"
Set-TextValue 48 3 "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue 48 4 "40.03"
Set-TextValue 48 5 "  -1.67%  "

Set-TextValue 49 2 "ONDO"
Set-TextValue 49 3 "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue 49 4 "1.31"
Set-TextValue 49 5 "  +11.21%  "

Set-TextValue 50 4 "141.36"
Set-TextValue 50 5 "  +0.14%  "

Set-TextValue 51 4 "2.787.61"
